$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 and 8 (the last two data rows), shrinking the used range to A1:I6
$ws.Rows.Item(7).Resize(2).Delete()

# Keep the Cpeid column (all "40105449000000E5") as text, not scientific-notation numbers
$ws.Range("A3:A6").NumberFormat = "@"

# Row 3: La Llorona - Le lacrime del male (voddashhttps CDN)
$ws.Cells.Item(3, 1).Value = "40105449000000E5"
$ws.Cells.Item(3, 2).Value = "cubo"
$ws.Cells.Item(3, 3).Value = "2019-08-13T07:37:21+0200"
$ws.Cells.Item(3, 4).Value = "2019-08-13T07:37:21+0200"
$ws.Cells.Item(3, 5).Value = "La Llorona - Le lacrime del male"
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = "voddashhttps.cb.ticdn.it"

# Row 4: Europei a squadre (livetv CDN), empty Start
$ws.Cells.Item(4, 1).Value = "40105449000000E5"
$ws.Cells.Item(4, 2).Value = "cubo"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = "2019-08-13T07:37:21+0200"
$ws.Cells.Item(4, 5).Value = "Europei a squadre"
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = "livetv0.cb.ticdn.it"

# Row 5: Dragon Trainer - Il mondo nascosto (voddashhttps CDN)
$ws.Cells.Item(5, 1).Value = "40105449000000E5"
$ws.Cells.Item(5, 2).Value = "cubo"
$ws.Cells.Item(5, 3).Value = "2019-08-14T08:30:04+0200"
$ws.Cells.Item(5, 4).Value = "2019-08-14T10:26:32+0200"
$ws.Cells.Item(5, 5).Value = "Dragon Trainer - Il mondo nascosto"
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = "voddashhttps.cb.ticdn.it"

# Row 6: Alita - Angelo della battaglia (voddashhttps CDN)
$ws.Cells.Item(6, 1).Value = "40105449000000E5"
$ws.Cells.Item(6, 2).Value = "cubo"
$ws.Cells.Item(6, 3).Value = "2019-08-13T07:37:21+0200"
$ws.Cells.Item(6, 4).Value = "2019-08-13T07:47:28+0200"
$ws.Cells.Item(6, 5).Value = "Alita - Angelo della battaglia"
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = "voddashhttps.cb.ticdn.it"
